# "fix(gui) step 1 and 2"
# Step 1: bump the list date in A1 by one day (2024-01-17 -> 2024-01-18,
# serial 45308 -> 45309).
# Step 2: update the price in D44 (Tender RODANTE Grande) to the new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45309
$ws.Range("D44").Value = 43783.243
